# B6-PowerPoint.pptx edit
#
# 1) Re-style the three tables (slides 14, 15, 16) from the custom
#    "Table_0" style to the built-in "{06378340-1D14-4BB6-B151-4134AE7AB329}"
#    table style.
# 2) Swap the presentation's applied theme colour scheme ("Integral" /
#    "Red Violet") for the plain "Office" colour scheme that previously
#    only lived in the unused theme1.xml part.

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------
$newTableStyle = "{06378340-1D14-4BB6-B151-4134AE7AB329}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newTableStyle)
    }
}

# --- 2) Theme colours --------------------------------------------------
function Hex-ToBgr($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office theme colour scheme, in MSO theme-colour-index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slideOne = $p.Slides.Item(1)
$themeColors = $slideOne.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = Hex-ToBgr $officeColors[$i - 1]
}
